$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'23.301.27"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.51%  "
$ws.Range("D3").Value = "'1.622.67"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.96%  "
$ws.Range("D4").Value = "'1.003"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.13%  "
$ws.Range("D5").Value = "'1.002"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.09%  "
$ws.Range("D6").Value = "'303.04"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.63%  "
$ws.Range("D7").Value = "'0.3731"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.26%  "
$ws.Range("B8").Value = "OKB"
$ws.Range("C8").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D8").Value = "'51.46"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.95%  "
$ws.Range("B9").Value = "Cardano"
$ws.Range("C9").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D9").Value = "'0.3613"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.22%  "
$ws.Range("E10").Value = "  +0.14%  "
$ws.Range("D11").Value = "'1.220"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.83%  "
$ws.Range("E12").Value = "  +0.04%  "
$ws.Range("D13").Value = "'22.19"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.90%  "
$ws.Range("D14").Value = "'6.454"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.24%  "
$ws.Range("D15").Value = "'0.00001236"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.54%  "
$ws.Range("D16").Value = "'7.263"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.31%  "
$ws.Range("D17").Value = "'1.621.76"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.98%  "
$ws.Range("D18").Value = "'93.71"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.57%  "
$ws.Range("D19").Value = "'0.06937"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.46%  "
$ws.Range("D20").Value = "'17.48"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -3.54%  "
$ws.Range("D21").Value = "'6.508"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.08%  "
$ws.Range("D22").Value = "'1.002"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.13%  "
$ws.Range("D23").Value = "'12.51"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.73%  "
$ws.Range("D24").Value = "'23.328.47"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.45%  "
$ws.Range("B25").Value = "Toncoin"
$ws.Range("C25").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D25").Value = "'2.457"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.55%  "
$ws.Range("B26").Value = "LidoDAOToken"
$ws.Range("C26").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D26").Value = "'3.093"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.06%  "
$ws.Range("D27").Value = "'21.12"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.38%  "
$ws.Range("D28").Value = "'149.68"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.16%  "
$ws.Range("D29").Value = "'5.237"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.50%  "
$ws.Range("D30").Value = "'132.44"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.46%  "
$ws.Range("D31").Value = "'1.803.14"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.80%  "
$ws.Range("D32").Value = "'6.694"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.60%  "
$ws.Range("D33").Value = "'2.120"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -6.67%  "
$ws.Range("D34").Value = "'1.042"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +9.14%  "
$ws.Range("D35").Value = "'10.80"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +4.82%  "
$ws.Range("D36").Value = "'0.02753"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.27%  "
$ws.Range("B37").Value = "Algorand"
$ws.Range("C37").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D37").Value = "'0.2491"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.03%  "
$ws.Range("B38").Value = "Stellar"
$ws.Range("C38").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D38").Value = "'0.08734"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.58%  "
$ws.Range("D39").Value = "'0.07076"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.38%  "
$ws.Range("D40").Value = "'5.959"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.79%  "
$ws.Range("D41").Value = "'0.6963"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.08%  "
$ws.Range("D42").Value = "'1.335"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.86%  "
$ws.Range("E43").Value = "  -0.37%  "
$ws.Range("D44").Value = "'12.03"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -3.27%  "
$ws.Range("D45").Value = "'0.6442"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.90%  "
$ws.Range("E46").Value = "  +0.07%  "
$ws.Range("E47").Value = "  -1.27%  "
$ws.Range("D48").Value = "'2.257"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.89%  "
$ws.Range("D49").Value = "'0.07960"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.24%  "
$ws.Range("D50").Value = "'125.40"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.23%  "
$ws.Range("D51").Value = "'1.180"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.70%  "
